# Add a new "Save" column (H) to the s_vals sheet.
# Column H mirrors the existing header styling (copied from G1) and is
# populated with a 0/1 indicator value for each data row (2-42).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1: text "Save", formatted like the other header cells.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for the new "Save" column, row by row.
$saveValues = @{
    2  = 1
    3  = 0
    4  = 1
    5  = 1
    6  = 0
    7  = 1
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 1
    14 = 0
    15 = 0
    16 = 1
    17 = 1
    18 = 0
    19 = 0
    20 = 1
    21 = 1
    22 = 0
    23 = 0
    24 = 0
    25 = 0
    26 = 0
    27 = 0
    28 = 1
    29 = 0
    30 = 0
    31 = 1
    32 = 0
    33 = 0
    34 = 1
    35 = 1
    36 = 1
    37 = 0
    38 = 1
    39 = 0
    40 = 0
    41 = 1
    42 = 1
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
